{"js": "// Locate the three paragraphs we need to touch by their current text:\n//  - Q28 paragraph: \"...CI m\u00ednimo (build + testes + Postgres + Flyway)?\"\n//  - the empty paragraph right after it (to be removed)\n//  - the paragraph that currently starts with \"29) Este projeto ...\"\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nlet q28Index = -1;\nlet emptyIndex = -1;\nlet q29Index = -1;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const t = paragraphs.items[i].text;\n  if (t.indexOf(\"CI m\u00ednimo (build + testes + Postgres + Flyway)?\") !== -1) {\n    q28Index = i;\n  } else if (q28Index !== -1 && emptyIndex === -1 && i === q28Index + 1 && t.trim() === \"\") {\n    emptyIndex = i;\n  } else if (t.indexOf(\"Este projeto vai se integrar com o mercado livre\") !== -1) {\n    q29Index = i;\n  }\n}\n\nconst q28Paragraph = paragraphs.items[q28Index];\nconst emptyParagraph = paragraphs.items[emptyIndex];\nconst q29Paragraph = paragraphs.items[q29Index];\n\n// 1) Remove the empty paragraph that used to separate question 28 and 29.\nemptyParagraph.delete();\n\n// 2) Append the new question 29 to the end of question 28's paragraph,\n//    separated by two manual line breaks (kept inside the same run/\n//    paragraph, matching a bold continuation of question 28's text).\nconst q28Tail = q28Paragraph.getRange(Word.RangeLocation.end);\nq28Tail.insertText(\n  \"\\u000b\\u000b29) Todo erro de regra usa ApiException com  ApiErrorCode correto?\",\n  Word.InsertLocation.replace\n);\n\n// 3) Insert a brand new bold paragraph for question 30 right before the\n//    \"Este projeto vai se integrar...\" paragraph (currently numbered 29).\nconst q30Paragraph = q29Paragraph.insertParagraph(\n  \"30) Me mostre quais os controllers que n\u00e3o est\u00e3o usando DTO\",\n  Word.InsertLocation.before\n);\nq30Paragraph.font.bold = true;\n// Also bold the paragraph mark itself (matches the authored formatting,\n// where the empty-paragraph mark carries bold/bCs run properties).\nq30Paragraph.getRange(Word.RangeLocation.end).font.bold = true;\n\n// 4) Renumber the \"29) Este projeto...\" paragraph to \"31) Este projeto...\".\nconst q29Scope = q29Paragraph.getRange();\nconst q29Marker = q29Scope.search(\"29)\", { matchCase: true });\nq29Marker.load(\"items\");\nawait context.sync();\nq29Marker.items[0].insertText(\"31)\", Word.InsertLocation.replace);\n\n// 5) Renumber the \"30) Caso tenha observado...\" question (inside the same\n//    paragraph as question 31) to \"32) Caso tenha observado...\".\nconst q32Scope = q29Paragraph.getRange();\nconst q32Marker = q32Scope.search(\"30)\", { matchCase: true });\nq32Marker.load(\"items\");\nawait context.sync();\nq32Marker.items[0].insertText(\"32)\", Word.InsertLocation.replace);\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# Locate the three paragraphs we need to touch by their current text:\n#  - Q28 paragraph: \"...CI m\u00ednimo (build + testes + Postgres + Flyway)?\"\n#  - the empty paragraph right after it (to be removed)\n#  - the paragraph that currently starts with \"29) Este projeto ...\"\n$count = $d.Paragraphs.Count\n$q28Index = -1\n$emptyIndex = -1\n$q29Index = -1\nfor ($i = 1; $i -le $count; $i++) {\n    $t = $d.Paragraphs.Item($i).Range.Text\n    if ($t.Contains(\"CI m\u00ednimo (build + testes + Postgres + Flyway)?\")) {\n        $q28Index = $i\n    } elseif ($q28Index -ne -1 -and $emptyIndex -eq -1 -and $i -eq ($q28Index + 1) -and $t.Trim().Length -eq 0) {\n        $emptyIndex = $i\n    } elseif ($t.Contains(\"Este projeto vai se integrar com o mercado livre\")) {\n        $q29Index = $i\n    }\n}\n\n# 1) Remove the empty paragraph that used to separate question 28 and 29.\n$d.Paragraphs.Item($emptyIndex).Range.Delete() | Out-Null\n\n# After the delete, question 29's paragraph shifted up by one.\n$q29Index = $q29Index - 1\n\n# 2) Append the new question 29 to the end of question 28's paragraph,\n#    separated by two manual line breaks (bold, matching question 28's\n#    run formatting). A Paragraph.Range's End sits *at* the paragraph\n#    mark, so shrink it by one character first or InsertAfter lands in\n#    the following paragraph instead of this one.\n$lineBreak = [char]11\n$q28Range = $d.Paragraphs.Item($q28Index).Range\n$q28Range.MoveEnd(1, -1) | Out-Null\n$q28Range.Collapse(0) | Out-Null\n$q28Range.InsertAfter(\"$lineBreak$lineBreak\" + \"29) Todo erro de regra usa ApiException com  ApiErrorCode correto?\")\n$q28Range.Font.Bold = $true\n$q28Range.Font.BoldBi = $true\n\n# 3) Insert a brand new bold paragraph for question 30 right before the\n#    \"Este projeto vai se integrar...\" paragraph (currently numbered 29).\n$insertionPoint = $d.Paragraphs.Item($q29Index).Range\n$insertionPoint.Collapse(1) | Out-Null\n$insertionPoint.InsertParagraphBefore() | Out-Null\n\n# The question-29 paragraph shifted down by one again; the freshly minted\n# empty paragraph now sits at the old $q29Index.\n$newParaIndex = $q29Index\n$q29Index = $q29Index + 1\n\n$newParaRange = $d.Paragraphs.Item($newParaIndex).Range\n$newParaRange.MoveEnd(1, -1) | Out-Null\n$newParaRange.InsertAfter(\"30) Me mostre quais os controllers que n\u00e3o est\u00e3o usando DTO\")\n$newParaRange.Font.Bold = $true\n$newParaRange.Font.BoldBi = $true\n\n# 4) Renumber the \"29) Este projeto...\" paragraph to \"31) Este projeto...\".\n$q29Range = $d.Paragraphs.Item($q29Index).Range\n$find29 = $q29Range.Find\n$find29.Execute(\"29)\") | Out-Null\n$q29Range.Text = \"31)\"\n\n# 5) Renumber the \"30) Caso tenha observado...\" question (inside the same\n#    paragraph as question 31) to \"32) Caso tenha observado...\".\n$q32Range = $d.Paragraphs.Item($q29Index).Range\n$find32 = $q32Range.Find\n$find32.Execute(\"30)\") | Out-Null\n$q32Range.Text = \"32)\"\n"}
